$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "90.796.73"
$ws.Range("E2").Value = "  +3.97%  "
$ws.Range("D3").Value = "3.213.95"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.18"
$ws.Range("E5").Value = "  +7.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "647.40"
$ws.Range("E6").Value = "  +6.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.405"
$ws.Range("E7").Value = "  +7.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.711"
$ws.Range("E8").Value = "  +7.59%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "3.209.24"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  +9.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.181"
$ws.Range("E12").Value = "  +2.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  +9.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.44"
$ws.Range("E14").Value = "  +3.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.62"
$ws.Range("E15").Value = "  +5.08%  "
$ws.Range("D16").Value = "90.456.42"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D17").Value = "3.792.40"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "3.203.12"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("E19").Value = "  +11.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000229"
$ws.Range("E20").Value = "  +78.00%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.59"
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "443.89"
$ws.Range("E22").Value = "  +7.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.71"
$ws.Range("E23").Value = "  +3.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.10"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.35"
$ws.Range("E25").Value = "  +3.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.05"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "82.26"
$ws.Range("E27").Value = "  +12.35%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.160"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.24"
$ws.Range("E32").Value = "  +43.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.50"
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "545.98"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.16"
$ws.Range("E35").Value = "  +6.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.94"
$ws.Range("E36").Value = "  +5.29%  "
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.66"
$ws.Range("E38").Value = "  +4.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.41"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.96"
$ws.Range("E42").Value = "  +3.45%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.376"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "147.17"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "44.92"
$ws.Range("E46").Value = "  +4.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "174.04"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.126"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.758"
$ws.Range("E49").Value = "  +9.14%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.627"
$ws.Range("E50").Value = "  +7.95%  "
$ws.Range("B51").Value = "ImmutableX"
$ws.Range("C51").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.25"
$ws.Range("E51").Value = "  +2.19%  "
